$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Тест дата" / "Тест время" test columns (O and P).
# Row 1 headers are removed entirely.
$ws.Range("O1").ClearContents()
$ws.Range("P1").ClearContents()

# Row 2 values are cleared but keep their existing number formatting/style.
$ws.Range("O2").ClearContents()
$ws.Range("P2").ClearContents()

# Update the view/selection to match the new state.
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("P2").Select()
